# Commit: "Fruta / hortaliza, semanal" — a new weekly price-report row is
# inserted into the daily log at row 98 (pushing all subsequent rows down
# by one), and populated with a new "Primera" quality Zanahoria record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 98; existing rows 98:204 shift down to 99:205.
$ws.Rows("98:98").Insert()

# Populate the newly inserted row with the reported data.
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 44494
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = 100114013
$ws.Range("G98").Value = "Zanahoria"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 13000
$ws.Range("L98").Value = 14000
$ws.Range("M98").Value = 13500
$ws.Range("N98").Value = "$/saco 25 kilos"
$ws.Range("O98").Value = "Valle de Camiña"
$ws.Range("P98").Value = 540
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"
